# Insert a new row at row 7 (weekly price update), pushing existing rows 7-9 down to 8-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

# Fill the newly inserted row 7 with the new week's data.
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value = 44907
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = 100112030
$ws.Cells.Item(7, 7).Value = "Poroto granado"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 2300
$ws.Cells.Item(7, 11).Value = 900
$ws.Cells.Item(7, 12).Value = 1000
$ws.Cells.Item(7, 13).Value = 952
$ws.Cells.Item(7, 14).Value = "$/kilo"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 952
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# Copy the date cell's style (s="2" -> custom date format) from the row below into the new row.
$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
